# Sanity semilla 8 - update "Cedula Cliente" / MSIDN sample values on the
# "Semilla 8" worksheet (sheet3.xml).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 8")

# Column B ("Cedula Cliente") values, rows 9-13 are stored as text.
$ws.Range("B9").Value  = "459399130"
$ws.Range("B10").Value = "836898669"
$ws.Range("B11").Value = "255188531"
$ws.Range("B12").Value = "194936717"
$ws.Range("B13").Value = "432694001"

# Column C ("MSIDN") values for rows 12-13.
$ws.Range("C12").Value = "3045981670"
$ws.Range("C13").Value = "3045981684"

# Row 14, column B is a real number (not text) in this sheet.
$ws.Range("B14").Value = 920626579

# Restore the selection/active cell Excel left the sheet on after editing.
$ws.Range("B10").Select()
